$d = $word.ActiveDocument

# Locate the "Edison Achalma" paragraph that uses the "Author" style
# (the byline right under the article title) robustly, rather than by a
# hard-coded paragraph index.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Style.NameLocal -eq "Author" -and $p.Range.Text.Trim() -eq "Edison Achalma") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not locate the 'Edison Achalma' Author paragraph"
}

$insertAt = $target.Range.End

# Build a range that spans the paragraph mark ending the "Edison Achalma"
# paragraph, then replace it (via InsertXML) with that same paragraph mark
# plus a brand-new "Author" styled paragraph containing the affiliation
# text. Doing it this way (rather than InsertParagraphAfter + Style=)
# keeps the original paragraph's OOXML untouched and avoids stray
# w:rsid*/w:paragraphStyle artifacts on the new paragraph.
$mark = $d.Range($insertAt - 1, $insertAt)

$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$xml = '<w:p ' + $ns + '><w:pPr><w:pStyle w:val="Author"/></w:pPr><w:r><w:t xml:space="preserve">Edison Achalma</w:t></w:r></w:p>' + `
       '<w:p ' + $ns + '><w:pPr><w:pStyle w:val="Author"/></w:pPr><w:r><w:t xml:space="preserve">Escuela Profesional de Economía, Universidad Nacional de San Cristóbal de Huamanga</w:t></w:r></w:p>'

[void]$mark.InsertXML($xml)

Write-Output "Inserted affiliation paragraph after 'Edison Achalma'"
